# Fruta / hortaliza, semanal
# Insert 4 new weekly rows at row 575 (pushing the existing rows 575-629 down to 579-633)
# and populate the new rows with the latest weekly price data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows at 575..578 (existing content shifts down)
$ws.Range("A575:A578").EntireRow.Insert()

# Common/static columns for all 4 new rows
$mercado = "Terminal Hortofrutícola Agro Chillán"
$region  = "Ñuble"
$codreg  = 16
$catId   = 100112004
$categoria = "Cebolla"
$clasificacion = "Hortaliza"

# Row 575
$r = 575
$ws.Cells.Item($r,1).Value2  = 7
$ws.Cells.Item($r,2).Value2  = $mercado
$ws.Cells.Item($r,3).Value2  = $region
$ws.Cells.Item($r,4).Value2  = 44769
$ws.Cells.Item($r,5).Value2  = $codreg
$ws.Cells.Item($r,6).Value2  = $catId
$ws.Cells.Item($r,7).Value2  = $categoria
$ws.Cells.Item($r,8).Value2  = "Sin especificar"
$ws.Cells.Item($r,9).Value2  = "1a (guarda)"
$ws.Cells.Item($r,10).Value2 = 120
$ws.Cells.Item($r,11).Value2 = 5000
$ws.Cells.Item($r,12).Value2 = 5500
$ws.Cells.Item($r,13).Value2 = 5250
$ws.Cells.Item($r,14).Value2 = "$/malla 15 kilos"
$ws.Cells.Item($r,15).Value2 = "Región del Maule"
$ws.Cells.Item($r,16).Value2 = 350
$ws.Cells.Item($r,17).Value2 = 15
$ws.Cells.Item($r,18).Value2 = $clasificacion

# Row 576
$r = 576
$ws.Cells.Item($r,1).Value2  = 7
$ws.Cells.Item($r,2).Value2  = $mercado
$ws.Cells.Item($r,3).Value2  = $region
$ws.Cells.Item($r,4).Value2  = 44769
$ws.Cells.Item($r,5).Value2  = $codreg
$ws.Cells.Item($r,6).Value2  = $catId
$ws.Cells.Item($r,7).Value2  = $categoria
$ws.Cells.Item($r,8).Value2  = "Sin especificar"
$ws.Cells.Item($r,9).Value2  = "1a (guarda)"
$ws.Cells.Item($r,10).Value2 = 120
$ws.Cells.Item($r,11).Value2 = 7000
$ws.Cells.Item($r,12).Value2 = 7500
$ws.Cells.Item($r,13).Value2 = 7250
$ws.Cells.Item($r,14).Value2 = "$/malla 25 kilos"
$ws.Cells.Item($r,15).Value2 = "Región del Maule"
$ws.Cells.Item($r,16).Value2 = 290
$ws.Cells.Item($r,17).Value2 = 25
$ws.Cells.Item($r,18).Value2 = $clasificacion

# Row 577
$r = 577
$ws.Cells.Item($r,1).Value2  = 7
$ws.Cells.Item($r,2).Value2  = $mercado
$ws.Cells.Item($r,3).Value2  = $region
$ws.Cells.Item($r,4).Value2  = 44769
$ws.Cells.Item($r,5).Value2  = $codreg
$ws.Cells.Item($r,6).Value2  = $catId
$ws.Cells.Item($r,7).Value2  = $categoria
$ws.Cells.Item($r,8).Value2  = "Sin especificar"
$ws.Cells.Item($r,9).Value2  = "2a (guarda)"
$ws.Cells.Item($r,10).Value2 = 80
$ws.Cells.Item($r,11).Value2 = 4500
$ws.Cells.Item($r,12).Value2 = 4500
$ws.Cells.Item($r,13).Value2 = 4500
$ws.Cells.Item($r,14).Value2 = "$/malla 15 kilos"
$ws.Cells.Item($r,15).Value2 = "Región del Maule"
$ws.Cells.Item($r,16).Value2 = 300
$ws.Cells.Item($r,17).Value2 = 15
$ws.Cells.Item($r,18).Value2 = $clasificacion

# Row 578
$r = 578
$ws.Cells.Item($r,1).Value2  = 7
$ws.Cells.Item($r,2).Value2  = $mercado
$ws.Cells.Item($r,3).Value2  = $region
$ws.Cells.Item($r,4).Value2  = 44769
$ws.Cells.Item($r,5).Value2  = $codreg
$ws.Cells.Item($r,6).Value2  = $catId
$ws.Cells.Item($r,7).Value2  = $categoria
$ws.Cells.Item($r,8).Value2  = "Sin especificar"
$ws.Cells.Item($r,9).Value2  = "2a (guarda)"
$ws.Cells.Item($r,10).Value2 = 120
$ws.Cells.Item($r,11).Value2 = 6000
$ws.Cells.Item($r,12).Value2 = 6500
$ws.Cells.Item($r,13).Value2 = 6250
$ws.Cells.Item($r,14).Value2 = "$/malla 25 kilos"
$ws.Cells.Item($r,15).Value2 = "Región del Maule"
$ws.Cells.Item($r,16).Value2 = 250
$ws.Cells.Item($r,17).Value2 = 25
$ws.Cells.Item($r,18).Value2 = $clasificacion
